$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 26.02.2022 01:15"

# Row 5 (Makro) price refresh, as produced by the AWS bash price-check script:
#  - B5 gets the newly observed price
#  - C5 keeps the previous price (what used to be in B5)
#  - D5 becomes a literal text delta string instead of a numeric difference
#  - E5 becomes a literal text timestamp string instead of a numeric Excel date
$ws.Range("B5").Value = 37.3
$ws.Range("C5").Value = 36.9

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "+0.4"
$ws.Range("D5").ClearFormats()

$ws.Range("E5").Value = "2022-02-26 01:17:15"
$ws.Range("E5").ClearFormats()
